$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 (from source row 25)
$ws.Range("F24").Value = "Botafogo SP"
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = "CRB"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2.2
$ws.Range("K24").Value = "24/04/2023 15:42"
$ws.Range("L24").Value = 2.43
$ws.Range("M24").Value = "29/04/2023 21:57"
$ws.Range("N24").Value = 3.09
$ws.Range("O24").Value = "24/04/2023 15:42"
$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = "29/04/2023 21:51"
$ws.Range("R24").Value = 3.54
$ws.Range("S24").Value = "24/04/2023 15:42"
$ws.Range("T24").Value = 3.46
$ws.Range("U24").Value = "29/04/2023 21:57"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-crb/jah69YO2/"

# Row 25 (from source row 24)
$ws.Range("F25").Value = "Chapecoense-SC"
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = "Ponte Preta"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 2.2
$ws.Range("K25").Value = "24/04/2023 15:42"
$ws.Range("L25").Value = 2.37
$ws.Range("M25").Value = "29/04/2023 21:59"
$ws.Range("N25").Value = 3.04
$ws.Range("O25").Value = "24/04/2023 15:42"
$ws.Range("P25").Value = 3.05
$ws.Range("Q25").Value = "29/04/2023 21:59"
$ws.Range("R25").Value = 3.84
$ws.Range("S25").Value = "24/04/2023 15:42"
$ws.Range("T25").Value = 3.52
$ws.Range("U25").Value = "29/04/2023 21:59"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-ponte-preta/vZ3M5jOR/"

# Row 31 (from source row 32)
$ws.Range("F31").Value = "CRB"
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = "Sampaio Correa"
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1.92
$ws.Range("K31").Value = "29/04/2023 22:13"
$ws.Range("L31").Value = 1.86
$ws.Range("M31").Value = "02/05/2023 23:59"
$ws.Range("N31").Value = 3.33
$ws.Range("O31").Value = "29/04/2023 22:13"
$ws.Range("P31").Value = 3.51
$ws.Range("Q31").Value = "02/05/2023 23:57"
$ws.Range("R31").Value = 4.47
$ws.Range("S31").Value = "29/04/2023 22:13"
$ws.Range("T31").Value = 4.68
$ws.Range("U31").Value = "02/05/2023 23:59"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/brazil/serie-b/crb-sampaio-correa/pf7UaLi6/"

# Row 32 (from source row 31)
$ws.Range("F32").Value = "Ponte Preta"
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = "Botafogo SP"
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2.15
$ws.Range("K32").Value = "29/04/2023 22:13"
$ws.Range("L32").Value = 2.03
$ws.Range("M32").Value = "02/05/2023 23:23"
$ws.Range("N32").Value = 3.06
$ws.Range("O32").Value = "29/04/2023 22:13"
$ws.Range("P32").Value = 3.25
$ws.Range("Q32").Value = "02/05/2023 23:23"
$ws.Range("R32").Value = 3.72
$ws.Range("S32").Value = "29/04/2023 22:13"
$ws.Range("T32").Value = 4.28
$ws.Range("U32").Value = "02/05/2023 23:23"
$ws.Range("V32").Value = "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-botafogo-sp/xAqcal1e/"

# Row 40 (from source row 41)
$ws.Range("F40").Value = "Sampaio Correa"
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = "Juventude"
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 2.3
$ws.Range("K40").Value = "04/05/2023 00:12"
$ws.Range("L40").Value = 2
$ws.Range("M40").Value = "06/05/2023 23:13"
$ws.Range("N40").Value = 3.08
$ws.Range("O40").Value = "04/05/2023 00:12"
$ws.Range("P40").Value = 3.44
$ws.Range("Q40").Value = "06/05/2023 23:13"
$ws.Range("R40").Value = 3.31
$ws.Range("S40").Value = "04/05/2023 00:12"
$ws.Range("T40").Value = 4.11
$ws.Range("U40").Value = "06/05/2023 23:13"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-esporte-clube-juventude/UNSaF0ig/"

# Row 41 (from source row 40)
$ws.Range("F41").Value = "Mirassol"
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = "Vila Nova FC"
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 2.07
$ws.Range("K41").Value = "06/05/2023 15:09"
$ws.Range("L41").Value = 1.95
$ws.Range("M41").Value = "06/05/2023 23:09"
$ws.Range("N41").Value = 3.06
$ws.Range("O41").Value = "06/05/2023 15:09"
$ws.Range("P41").Value = 3.3
$ws.Range("Q41").Value = "06/05/2023 23:10"
$ws.Range("R41").Value = 3.96
$ws.Range("S41").Value = "06/05/2023 15:09"
$ws.Range("T41").Value = 4.55
$ws.Range("U41").Value = "06/05/2023 23:14"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/brazil/serie-b/mirassol-vila-nova-fc/YeZjHMMt/"

# Row 43 (from source row 44)
$ws.Range("F43").Value = "Chapecoense-SC"
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = "Novorizontino"
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = 2.21
$ws.Range("K43").Value = "04/05/2023 00:12"
$ws.Range("L43").Value = 2.54
$ws.Range("M43").Value = "07/05/2023 01:23"
$ws.Range("N43").Value = 3.13
$ws.Range("O43").Value = "04/05/2023 00:12"
$ws.Range("P43").Value = 3.05
$ws.Range("Q43").Value = "07/05/2023 01:23"
$ws.Range("R43").Value = 3.47
$ws.Range("S43").Value = "04/05/2023 00:12"
$ws.Range("T43").Value = 3.21
$ws.Range("U43").Value = "07/05/2023 01:23"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-novorizontino/OpS3EK6a/"

# Row 44 (from source row 43)
$ws.Range("F44").Value = "CRB"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "ABC"
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 1.89
$ws.Range("K44").Value = "04/05/2023 00:12"
$ws.Range("L44").Value = 1.72
$ws.Range("M44").Value = "07/05/2023 01:29"
$ws.Range("N44").Value = 3.24
$ws.Range("O44").Value = "04/05/2023 00:12"
$ws.Range("P44").Value = 3.52
$ws.Range("Q44").Value = "07/05/2023 01:29"
$ws.Range("R44").Value = 4.42
$ws.Range("S44").Value = "04/05/2023 00:12"
$ws.Range("T44").Value = 5.77
$ws.Range("U44").Value = "07/05/2023 01:29"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/brazil/serie-b/crb-abc/KCLCCbyC/"

# Row 60 (from source row 61)
$ws.Range("F60").Value = "Tombense"
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = "Londrina"
$ws.Range("I60").Value = 2
$ws.Range("J60").Value = 2
$ws.Range("K60").Value = "14/05/2023 20:42"
$ws.Range("L60").Value = 2.07
$ws.Range("M60").Value = "19/05/2023 23:58"
$ws.Range("N60").Value = 3.16
$ws.Range("O60").Value = "14/05/2023 20:42"
$ws.Range("P60").Value = 3.27
$ws.Range("Q60").Value = "19/05/2023 23:58"
$ws.Range("R60").Value = 4.08
$ws.Range("S60").Value = "14/05/2023 20:42"
$ws.Range("T60").Value = 4.06
$ws.Range("U60").Value = "19/05/2023 23:58"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/brazil/serie-b/tombense-londrina/SK5UM3GG/"

# Row 61 (from source row 60)
$ws.Range("F61").Value = "Mirassol"
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = "Vitoria"
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2.28
$ws.Range("K61").Value = "14/05/2023 23:12"
$ws.Range("L61").Value = 2.22
$ws.Range("M61").Value = "19/05/2023 23:52"
$ws.Range("N61").Value = 3.05
$ws.Range("O61").Value = "14/05/2023 23:12"
$ws.Range("P61").Value = 3.22
$ws.Range("Q61").Value = "19/05/2023 23:52"
$ws.Range("R61").Value = 3.4
$ws.Range("S61").Value = "14/05/2023 23:12"
$ws.Range("T61").Value = 3.68
$ws.Range("U61").Value = "19/05/2023 23:52"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/brazil/serie-b/mirassol-vitoria/buEDQoHi/"

# Row 64 (from source row 65)
$ws.Range("F64").Value = "Sampaio Correa"
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = "ABC"
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1.93
$ws.Range("K64").Value = "14/05/2023 23:12"
$ws.Range("L64").Value = 1.91
$ws.Range("M64").Value = "20/05/2023 21:52"
$ws.Range("N64").Value = 3.34
$ws.Range("O64").Value = "14/05/2023 23:12"
$ws.Range("P64").Value = 3.47
$ws.Range("Q64").Value = "20/05/2023 21:57"
$ws.Range("R64").Value = 4.4
$ws.Range("S64").Value = "14/05/2023 23:12"
$ws.Range("T64").Value = 4.47
$ws.Range("U64").Value = "20/05/2023 21:57"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-abc/f77MOPo4/"

# Row 65 (from source row 64)
$ws.Range("F65").Value = "Ituano"
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = "Novorizontino"
$ws.Range("I65").Value = 2
$ws.Range("J65").Value = 2.65
$ws.Range("K65").Value = "15/05/2023 01:42"
$ws.Range("L65").Value = 2.71
$ws.Range("M65").Value = "20/05/2023 21:59"
$ws.Range("N65").Value = 3.13
$ws.Range("O65").Value = "15/05/2023 01:42"
$ws.Range("P65").Value = 3.08
$ws.Range("Q65").Value = "20/05/2023 21:59"
$ws.Range("R65").Value = 2.75
$ws.Range("S65").Value = "15/05/2023 01:42"
$ws.Range("T65").Value = 2.94
$ws.Range("U65").Value = "20/05/2023 21:50"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/brazil/serie-b/ituano-novorizontino/KlF9RR1o/"

# Row 73 (from source row 76)
$ws.Range("F73").Value = "Vitoria"
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = "CRB"
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1.78
$ws.Range("K73").Value = "21/05/2023 23:42"
$ws.Range("L73").Value = 1.71
$ws.Range("M73").Value = "24/05/2023 23:37"
$ws.Range("N73").Value = 3.35
$ws.Range("O73").Value = "21/05/2023 23:42"
$ws.Range("P73").Value = 3.41
$ws.Range("Q73").Value = "24/05/2023 23:59"
$ws.Range("R73").Value = 4.92
$ws.Range("S73").Value = "21/05/2023 23:42"
$ws.Range("T73").Value = 5.69
$ws.Range("U73").Value = "24/05/2023 23:59"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/brazil/serie-b/vitoria-crb/SzfhxqG3/"

# Row 76 (from source row 73)
$ws.Range("F76").Value = "Juventude"
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = "Atletico GO"
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 2.5
$ws.Range("K76").Value = "21/05/2023 16:12"
$ws.Range("L76").Value = 2.15
$ws.Range("M76").Value = "24/05/2023 23:59"
$ws.Range("N76").Value = 3.12
$ws.Range("O76").Value = "21/05/2023 16:12"
$ws.Range("P76").Value = 3.35
$ws.Range("Q76").Value = "24/05/2023 23:59"
$ws.Range("R76").Value = 3.1
$ws.Range("S76").Value = "21/05/2023 16:12"
$ws.Range("T76").Value = 3.7
$ws.Range("U76").Value = "24/05/2023 23:59"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/brazil/serie-b/esporte-clube-juventude-atletico-go/ljr1zNpG/"

# Row 80 (from source row 81)
$ws.Range("F80").Value = "CRB"
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = "Juventude"
$ws.Range("I80").Value = 2
$ws.Range("J80").Value = 2.12
$ws.Range("K80").Value = "25/05/2023 02:12"
$ws.Range("L80").Value = 2.53
$ws.Range("M80").Value = "27/05/2023 21:52"
$ws.Range("N80").Value = 3.22
$ws.Range("O80").Value = "25/05/2023 02:12"
$ws.Range("P80").Value = 3.08
$ws.Range("Q80").Value = "27/05/2023 21:52"
$ws.Range("R80").Value = 3.57
$ws.Range("S80").Value = "25/05/2023 02:12"
$ws.Range("T80").Value = 3.2
$ws.Range("U80").Value = "27/05/2023 21:52"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/brazil/serie-b/crb-esporte-clube-juventude/f9YFUiEf/"

# Row 81 (from source row 80)
$ws.Range("F81").Value = "Ituano"
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = "Londrina"
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 2.08
$ws.Range("K81").Value = "25/05/2023 02:12"
$ws.Range("L81").Value = 1.97
$ws.Range("M81").Value = "27/05/2023 21:31"
$ws.Range("N81").Value = 3.14
$ws.Range("O81").Value = "25/05/2023 02:12"
$ws.Range("P81").Value = 3.24
$ws.Range("Q81").Value = "27/05/2023 21:31"
$ws.Range("R81").Value = 3.8
$ws.Range("S81").Value = "25/05/2023 02:12"
$ws.Range("T81").Value = 4.56
$ws.Range("U81").Value = "27/05/2023 22:00"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/brazil/serie-b/ituano-londrina/Sh5Swdi8/"

# Row 92 (from source row 93)
$ws.Range("F92").Value = "Criciuma"
$ws.Range("G92").Value = 3
$ws.Range("H92").Value = "Atletico GO"
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 2.06
$ws.Range("K92").Value = "31/05/2023 05:12"
$ws.Range("L92").Value = 2.11
$ws.Range("M92").Value = "03/06/2023 02:13"
$ws.Range("N92").Value = 3.24
$ws.Range("O92").Value = "31/05/2023 05:12"
$ws.Range("P92").Value = 3.1
$ws.Range("Q92").Value = "03/06/2023 02:13"
$ws.Range("R92").Value = 3.73
$ws.Range("S92").Value = "31/05/2023 05:12"
$ws.Range("T92").Value = 4.21
$ws.Range("U92").Value = "03/06/2023 02:13"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/brazil/serie-b/criciuma-atletico-go/tQzSRkbD/"

# Row 93 (from source row 92)
$ws.Range("F93").Value = "Vitoria"
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = "Ituano"
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1.69
$ws.Range("K93").Value = "31/05/2023 05:12"
$ws.Range("L93").Value = 1.78
$ws.Range("M93").Value = "03/06/2023 02:24"
$ws.Range("N93").Value = 3.53
$ws.Range("O93").Value = "31/05/2023 05:12"
$ws.Range("P93").Value = 3.51
$ws.Range("Q93").Value = "03/06/2023 02:24"
$ws.Range("R93").Value = 5.21
$ws.Range("S93").Value = "31/05/2023 05:12"
$ws.Range("T93").Value = 5.29
$ws.Range("U93").Value = "03/06/2023 02:24"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/brazil/serie-b/vitoria-ituano/jeIIlD5Q/"

# Row 100 (from source row 101)
$ws.Range("F100").Value = "Tombense"
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = "Vitoria"
$ws.Range("I100").Value = 2
$ws.Range("J100").Value = 2.46
$ws.Range("K100").Value = "03/06/2023 02:42"
$ws.Range("L100").Value = 2.46
$ws.Range("M100").Value = "06/06/2023 23:59"
$ws.Range("N100").Value = 3.03
$ws.Range("O100").Value = "03/06/2023 02:42"
$ws.Range("P100").Value = 3.18
$ws.Range("Q100").Value = "06/06/2023 23:59"
$ws.Range("R100").Value = 3.08
$ws.Range("S100").Value = "03/06/2023 02:42"
$ws.Range("T100").Value = 3.2
$ws.Range("U100").Value = "06/06/2023 23:59"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/brazil/serie-b/tombense-vitoria/Kp3SflSa/"

# Row 101 (from source row 102)
$ws.Range("F101").Value = "Ponte Preta"
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = "ABC"
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 1.7
$ws.Range("K101").Value = "03/06/2023 22:12"
$ws.Range("L101").Value = 1.58
$ws.Range("M101").Value = "06/06/2023 23:38"
$ws.Range("N101").Value = 3.56
$ws.Range("O101").Value = "03/06/2023 22:12"
$ws.Range("P101").Value = 3.74
$ws.Range("Q101").Value = "06/06/2023 23:38"
$ws.Range("R101").Value = 5.08
$ws.Range("S101").Value = "03/06/2023 22:12"
$ws.Range("T101").Value = 7.21
$ws.Range("U101").Value = "06/06/2023 23:38"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-abc/fT01tkD6/"

# Row 102 (from source row 100)
$ws.Range("F102").Value = "Chapecoense-SC"
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = "Vila Nova FC"
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = 2.37
$ws.Range("K102").Value = "04/06/2023 00:13"
$ws.Range("L102").Value = 2.81
$ws.Range("M102").Value = "06/06/2023 23:59"
$ws.Range("N102").Value = 3.05
$ws.Range("O102").Value = "04/06/2023 00:13"
$ws.Range("P102").Value = 2.99
$ws.Range("Q102").Value = "06/06/2023 23:59"
$ws.Range("R102").Value = 3.22
$ws.Range("S102").Value = "04/06/2023 00:13"
$ws.Range("T102").Value = 2.92
$ws.Range("U102").Value = "06/06/2023 23:59"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-vila-nova-fc/6saOeUCg/"

# Row 163 (from source row 164)
$ws.Range("F163").Value = "Atletico GO"
$ws.Range("G163").Value = 3
$ws.Range("H163").Value = "Sport Recife"
$ws.Range("I163").Value = 1
$ws.Range("J163").Value = 2.64
$ws.Range("K163").Value = "09/07/2023 23:11"
$ws.Range("L163").Value = 2.68
$ws.Range("M163").Value = "15/07/2023 02:27"
$ws.Range("N163").Value = 3.01
$ws.Range("O163").Value = "09/07/2023 23:11"
$ws.Range("P163").Value = 2.99
$ws.Range("Q163").Value = "15/07/2023 02:27"
$ws.Range("R163").Value = 2.87
$ws.Range("S163").Value = "09/07/2023 23:11"
$ws.Range("T163").Value = 3.07
$ws.Range("U163").Value = "15/07/2023 02:29"
$ws.Range("V163").Value = "https://www.betexplorer.com/football/brazil/serie-b/atletico-go-sport-recife/KlNQMIY6/"

# Row 164 (from source row 163)
$ws.Range("F164").Value = "Sampaio Correa"
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = "Ituano"
$ws.Range("I164").Value = 0
$ws.Range("J164").Value = 1.99
$ws.Range("K164").Value = "08/07/2023 23:12"
$ws.Range("L164").Value = 1.95
$ws.Range("M164").Value = "15/07/2023 02:21"
$ws.Range("N164").Value = 3.21
$ws.Range("O164").Value = "08/07/2023 23:12"
$ws.Range("P164").Value = 3.2
$ws.Range("Q164").Value = "15/07/2023 02:20"
$ws.Range("R164").Value = 4.35
$ws.Range("S164").Value = "08/07/2023 23:12"
$ws.Range("T164").Value = 4.74
$ws.Range("U164").Value = "15/07/2023 02:21"
$ws.Range("V164").Value = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-ituano/KnWlRXPh/"

# Row 176 (from source row 178)
$ws.Range("F176").Value = "ABC"
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = "Guarani"
$ws.Range("I176").Value = 1
$ws.Range("J176").Value = 2.6
$ws.Range("K176").Value = "15/07/2023 16:12"
$ws.Range("L176").Value = 2.6
$ws.Range("M176").Value = "20/07/2023 02:29"
$ws.Range("N176").Value = 2.85
$ws.Range("O176").Value = "15/07/2023 16:12"
$ws.Range("P176").Value = 2.81
$ws.Range("Q176").Value = "20/07/2023 02:26"
$ws.Range("R176").Value = 3.08
$ws.Range("S176").Value = "15/07/2023 16:12"
$ws.Range("T176").Value = 3.42
$ws.Range("U176").Value = "20/07/2023 02:29"
$ws.Range("V176").Value = "https://www.betexplorer.com/football/brazil/serie-b/abc-guarani/KS9yy9gH/"

# Row 178 (from source row 176)
$ws.Range("F178").Value = "Sport Recife"
$ws.Range("G178").Value = 1
$ws.Range("H178").Value = "Vitoria"
$ws.Range("I178").Value = 2
$ws.Range("J178").Value = 1.69
$ws.Range("K178").Value = "16/07/2023 23:12"
$ws.Range("L178").Value = 1.74
$ws.Range("M178").Value = "20/07/2023 01:57"
$ws.Range("N178").Value = 3.51
$ws.Range("O178").Value = "16/07/2023 23:12"
$ws.Range("P178").Value = 3.42
$ws.Range("Q178").Value = "20/07/2023 01:54"
$ws.Range("R178").Value = 5.8
$ws.Range("S178").Value = "16/07/2023 23:12"
$ws.Range("T178").Value = 5.96
$ws.Range("U178").Value = "20/07/2023 01:57"
$ws.Range("V178").Value = "https://www.betexplorer.com/football/brazil/serie-b/sport-recife-vitoria/SCeSvifh/"

# Row 194 (from source row 195)
$ws.Range("F194").Value = "ABC"
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = "Londrina"
$ws.Range("I194").Value = 0
$ws.Range("J194").Value = 2.05
$ws.Range("K194").Value = "25/07/2023 00:12"
$ws.Range("L194").Value = 2.18
$ws.Range("M194").Value = "29/07/2023 21:51"
$ws.Range("N194").Value = 3.05
$ws.Range("O194").Value = "25/07/2023 00:12"
$ws.Range("P194").Value = 3.07
$ws.Range("Q194").Value = "29/07/2023 21:51"
$ws.Range("R194").Value = 4.35
$ws.Range("S194").Value = "25/07/2023 00:12"
$ws.Range("T194").Value = 4.02
$ws.Range("U194").Value = "29/07/2023 21:51"
$ws.Range("V194").Value = "https://www.betexplorer.com/football/brazil/serie-b/abc-londrina/rq200PMS/"

# Row 195 (from source row 194)
$ws.Range("F195").Value = "Avai"
$ws.Range("G195").Value = 1
$ws.Range("H195").Value = "Guarani"
$ws.Range("I195").Value = 2
$ws.Range("J195").Value = 2.23
$ws.Range("K195").Value = "23/07/2023 16:12"
$ws.Range("L195").Value = 2.4
$ws.Range("M195").Value = "29/07/2023 21:58"
$ws.Range("N195").Value = 3.02
$ws.Range("O195").Value = "23/07/2023 16:12"
$ws.Range("P195").Value = 2.86
$ws.Range("Q195").Value = "29/07/2023 21:57"
$ws.Range("R195").Value = 3.78
$ws.Range("S195").Value = "23/07/2023 16:12"
$ws.Range("T195").Value = 3.74
$ws.Range("U195").Value = "29/07/2023 21:58"
$ws.Range("V195").Value = "https://www.betexplorer.com/football/brazil/serie-b/avai-guarani/QsYhspGj/"

# Row 200 (from source row 201)
$ws.Range("F200").Value = "Mirassol"
$ws.Range("G200").Value = 2
$ws.Range("H200").Value = "Avai"
$ws.Range("I200").Value = 2
$ws.Range("J200").Value = 1.57
$ws.Range("K200").Value = "29/07/2023 22:12"
$ws.Range("L200").Value = 1.66
$ws.Range("M200").Value = "01/08/2023 23:33"
$ws.Range("N200").Value = 3.79
$ws.Range("O200").Value = "29/07/2023 22:12"
$ws.Range("P200").Value = 3.68
$ws.Range("Q200").Value = "01/08/2023 23:51"
$ws.Range("R200").Value = 6.77
$ws.Range("S200").Value = "29/07/2023 22:12"
$ws.Range("T200").Value = 6.16
$ws.Range("U200").Value = "01/08/2023 23:33"
$ws.Range("V200").Value = "https://www.betexplorer.com/football/brazil/serie-b/mirassol-avai/6cZMHQi3/"

# Row 201 (from source row 200)
$ws.Range("F201").Value = "Sampaio Correa"
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = "Botafogo SP"
$ws.Range("I201").Value = 0
$ws.Range("J201").Value = 2.02
$ws.Range("K201").Value = "29/07/2023 23:12"
$ws.Range("L201").Value = 2.09
$ws.Range("M201").Value = "01/08/2023 23:51"
$ws.Range("N201").Value = 3.06
$ws.Range("O201").Value = "29/07/2023 23:12"
$ws.Range("P201").Value = 3.06
$ws.Range("Q201").Value = "01/08/2023 23:51"
$ws.Range("R201").Value = 4.16
$ws.Range("S201").Value = "29/07/2023 23:12"
$ws.Range("T201").Value = 4.39
$ws.Range("U201").Value = "01/08/2023 23:51"
$ws.Range("V201").Value = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-botafogo-sp/SCNVF4MF/"

# Row 205 (from source row 209)
$ws.Range("F205").Value = "Vila Nova FC"
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = "Sport Recife"
$ws.Range("I205").Value = 1
$ws.Range("J205").Value = 2.3
$ws.Range("K205").Value = "30/07/2023 20:42"
$ws.Range("L205").Value = 2.31
$ws.Range("M205").Value = "03/08/2023 02:19"
$ws.Range("N205").Value = 2.97
$ws.Range("O205").Value = "30/07/2023 20:42"
$ws.Range("P205").Value = 2.9
$ws.Range("Q205").Value = "03/08/2023 02:19"
$ws.Range("R205").Value = 3.68
$ws.Range("S205").Value = "30/07/2023 20:42"
$ws.Range("T205").Value = 3.9
$ws.Range("U205").Value = "03/08/2023 02:27"
$ws.Range("V205").Value = "https://www.betexplorer.com/football/brazil/serie-b/vila-nova-fc-sport-recife/2TrXf3pj/"

# Row 206 (from source row 205)
$ws.Range("F206").Value = "Londrina"
$ws.Range("G206").Value = 1
$ws.Range("H206").Value = "Chapecoense-SC"
$ws.Range("I206").Value = 1
$ws.Range("J206").Value = 2.51
$ws.Range("K206").Value = "29/07/2023 22:12"
$ws.Range("L206").Value = 2.59
$ws.Range("M206").Value = "03/08/2023 02:18"
$ws.Range("N206").Value = 2.96
$ws.Range("O206").Value = "29/07/2023 22:12"
$ws.Range("P206").Value = 2.9
$ws.Range("Q206").Value = "03/08/2023 02:18"
$ws.Range("R206").Value = 3.26
$ws.Range("S206").Value = "29/07/2023 22:12"
$ws.Range("T206").Value = 3.31
$ws.Range("U206").Value = "03/08/2023 02:18"
$ws.Range("V206").Value = "https://www.betexplorer.com/football/brazil/serie-b/londrina-chapecoense-sc/YBk3cd6l/"

# Row 207 (from source row 206)
$ws.Range("F207").Value = "Guarani"
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = "Ceara"
$ws.Range("I207").Value = 0
$ws.Range("J207").Value = 2.2
$ws.Range("K207").Value = "29/07/2023 22:12"
$ws.Range("L207").Value = 2.1
$ws.Range("M207").Value = "03/08/2023 02:28"
$ws.Range("N207").Value = 3.08
$ws.Range("O207").Value = "29/07/2023 22:12"
$ws.Range("P207").Value = 3.18
$ws.Range("Q207").Value = "03/08/2023 02:24"
$ws.Range("R207").Value = 3.79
$ws.Range("S207").Value = "29/07/2023 22:12"
$ws.Range("T207").Value = 4.11
$ws.Range("U207").Value = "03/08/2023 02:28"
$ws.Range("V207").Value = "https://www.betexplorer.com/football/brazil/serie-b/guarani-ceara/phlSeqVq/"

# Row 208 (from source row 207)
$ws.Range("F208").Value = "Juventude"
$ws.Range("G208").Value = 1
$ws.Range("H208").Value = "Novorizontino"
$ws.Range("I208").Value = 0
$ws.Range("J208").Value = 2.38
$ws.Range("K208").Value = "30/07/2023 20:42"
$ws.Range("L208").Value = 2.46
$ws.Range("M208").Value = "03/08/2023 02:21"
$ws.Range("N208").Value = 2.96
$ws.Range("O208").Value = "30/07/2023 20:42"
$ws.Range("P208").Value = 2.94
$ws.Range("Q208").Value = "03/08/2023 02:28"
$ws.Range("R208").Value = 3.51
$ws.Range("S208").Value = "30/07/2023 20:42"
$ws.Range("T208").Value = 3.48
$ws.Range("U208").Value = "03/08/2023 02:21"
$ws.Range("V208").Value = "https://www.betexplorer.com/football/brazil/serie-b/esporte-clube-juventude-novorizontino/fLjabxir/"

# Row 209 (from source row 208)
$ws.Range("F209").Value = "Vitoria"
$ws.Range("G209").Value = 2
$ws.Range("H209").Value = "ABC"
$ws.Range("I209").Value = 0
$ws.Range("J209").Value = 1.6
$ws.Range("K209").Value = "30/07/2023 23:12"
$ws.Range("L209").Value = 1.47
$ws.Range("M209").Value = "03/08/2023 02:29"
$ws.Range("N209").Value = 3.78
$ws.Range("O209").Value = "30/07/2023 23:12"
$ws.Range("P209").Value = 4.23
$ws.Range("Q209").Value = "03/08/2023 02:29"
$ws.Range("R209").Value = 5.67
$ws.Range("S209").Value = "30/07/2023 23:12"
$ws.Range("T209").Value = 8.199999999999999
$ws.Range("U209").Value = "03/08/2023 02:29"
$ws.Range("V209").Value = "https://www.betexplorer.com/football/brazil/serie-b/vitoria-abc/lMORGp79/"

# Row 235 (from source row 236)
$ws.Range("F235").Value = "Ituano"
$ws.Range("G235").Value = 3
$ws.Range("H235").Value = "Criciuma"
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = 2.69
$ws.Range("K235").Value = "15/08/2023 00:12"
$ws.Range("L235").Value = 2.79
$ws.Range("M235").Value = "19/08/2023 21:51"
$ws.Range("N235").Value = 3
$ws.Range("O235").Value = "15/08/2023 00:12"
$ws.Range("P235").Value = 2.95
$ws.Range("Q235").Value = "19/08/2023 21:51"
$ws.Range("R235").Value = 2.83
$ws.Range("S235").Value = "15/08/2023 00:12"
$ws.Range("T235").Value = 2.98
$ws.Range("U235").Value = "19/08/2023 21:51"
$ws.Range("V235").Value = "https://www.betexplorer.com/football/brazil/serie-b/ituano-criciuma/jmzKsEBP/"

# Row 236 (from source row 235)
$ws.Range("F236").Value = "ABC"
$ws.Range("G236").Value = 1
$ws.Range("H236").Value = "CRB"
$ws.Range("I236").Value = 2
$ws.Range("J236").Value = 2.62
$ws.Range("K236").Value = "15/08/2023 00:12"
$ws.Range("L236").Value = 2.94
$ws.Range("M236").Value = "19/08/2023 21:51"
$ws.Range("N236").Value = 2.95
$ws.Range("O236").Value = "15/08/2023 00:12"
$ws.Range("P236").Value = 2.8
$ws.Range("Q236").Value = "19/08/2023 21:57"
$ws.Range("R236").Value = 2.95
$ws.Range("S236").Value = "15/08/2023 00:12"
$ws.Range("T236").Value = 2.99
$ws.Range("U236").Value = "19/08/2023 21:51"
$ws.Range("V236").Value = "https://www.betexplorer.com/football/brazil/serie-b/abc-crb/jBI6Y9YB/"

# New row 354
$ws.Range("A354").Value = 353
$ws.Range("B354").Value = "brazil"
$ws.Range("C354").Value = "serie-b"
$ws.Range("D354").Value = "2023"
$ws.Range("E354").Value = 45241.8125
$ws.Range("F354").Value = "Ituano"
$ws.Range("G354").Value = 2
$ws.Range("H354").Value = "Sampaio Correa"
$ws.Range("I354").Value = 1
$ws.Range("J354").Value = 1.98
$ws.Range("K354").Value = "04/11/2023 23:43"
$ws.Range("L354").Value = 1.86
$ws.Range("M354").Value = "11/11/2023 19:22"
$ws.Range("N354").Value = 3.12
$ws.Range("O354").Value = "04/11/2023 23:43"
$ws.Range("P354").Value = 3.18
$ws.Range("Q354").Value = "11/11/2023 19:22"
$ws.Range("R354").Value = 4.58
$ws.Range("S354").Value = "04/11/2023 23:43"
$ws.Range("T354").Value = 5.42
$ws.Range("U354").Value = "11/11/2023 19:22"
$ws.Range("V354").Value = "https://www.betexplorer.com/football/brazil/serie-b/ituano-sampaio-correa/0hr9gopP/"

# Copy style/format for A354 and E354 from row 353 equivalents
$ws.Range("A353").Copy()
$ws.Range("A354").PasteSpecial(-4122)
$ws.Range("E353").Copy()
$ws.Range("E354").PasteSpecial(-4122)

Write-Output "Edit complete"